$wb = $excel.ActiveWorkbook

# --- BCbVT-passenger & BCbVT-freight: add "LPG vehicle" and "hydrogen vehicle"
# columns (inserted just before the existing "nonroad vehicle" column), fill
# the new columns with 0 for every vehicle-type row, label the row header in
# A1, and wrap/resize the header row.
foreach ($name in @("BCbVT-passenger", "BCbVT-freight")) {
    $ws = $wb.Worksheets.Item($name)

    # Insert two new columns at G:H (the former "nonroad vehicle" column,
    # G, shifts right to I).
    $ws.Range("G1:H1").EntireColumn.Insert()

    $ws.Range("G1").Value = "LPG vehicle"
    $ws.Range("H1").Value = "hydrogen vehicle"

    $ws.Range("G2:H7").Value = 0

    # Row label for the header row.
    $ws.Range("A1").Value = "Battery Capacity (MW*hr/vehicle"
    $ws.Range("A1").WrapText = $true

    $ws.Rows.Item(1).RowHeight = 57
}
